$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("M2").Value = 0.6327629999999999
$ws.Range("N2").Value = 1.898289
$ws.Range("O2").Value = 0.1382544270550543
$ws.Range("P2").Value = 0.1382544270550544
$ws.Range("Q2").Value = 107.516340870291
$ws.Range("R2").Value = 967.6470678326189
$ws.Range("S2").Value = 0.06140134781824951
$ws.Range("T2").Value = 0.06140134781824953
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("O3").Value = 0.4765301499162115
$ws.Range("P3").Value = 0.4765301499162115
$ws.Range("Q3").Value = 370.5832726279357
$ws.Range("R3").Value = 3335.249453651421
$ws.Range("S3").Value = 0.2116358521325065
$ws.Range("T3").Value = 0.2116358521325065
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("M4").Value = 1.444396333333334
$ws.Range("N4").Value = 4.333189000000001
$ws.Range("O4").Value = 0.3155908096798033
$ws.Range("P4").Value = 0.3155908096798033
$ws.Range("Q4").Value = 245.4255519467244
$ws.Range("R4").Value = 2208.829967520519
$ws.Range("S4").Value = 0.1401597148543836
$ws.Range("T4").Value = 0.1401597148543836
$ws.Range("G5").Value = 169.915657
$ws.Range("H5").Value = 509.746971
$ws.Range("I5").Value = 0.4441184931734509
$ws.Range("J5").Value = 0.4441184931734509
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3186579999999999
$ws.Range("N5").Value = 0.9559739999999999
$ws.Range("O5").Value = 0.06962461334893082
$ws.Range("P5").Value = 0.06962461334893082
$ws.Range("Q5").Value = 54.14498342830599
$ws.Range("R5").Value = 487.304850854754
$ws.Range("S5").Value = 0.03092157836831129
$ws.Range("T5").Value = 0.03092157836831129
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("M6").Value = 0.6327629999999999
$ws.Range("N6").Value = 1.898289
$ws.Range("O6").Value = 0.1382544270550543
$ws.Range("P6").Value = 0.1382544270550544
$ws.Range("Q6").Value = 43.26975322740899
$ws.Range("R6").Value = 389.4277790466809
$ws.Range("S6").Value = 0.02471085926492965
$ws.Range("T6").Value = 0.02471085926492965
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("O7").Value = 0.4765301499162115
$ws.Range("P7").Value = 0.4765301499162115
$ws.Range("S7").Value = 0.08517245863950682
$ws.Range("T7").Value = 0.08517245863950683
$ws.Range("I8").Value = 0.1787346690539575
$ws.Range("J8").Value = 0.1787346690539575
$ws.Range("M8").Value = 1.444396333333334
$ws.Range("N8").Value = 4.333189000000001
$ws.Range("O8").Value = 0.3155908096798033
$ws.Range("P8").Value = 0.3155908096798033
$ws.Range("Q8").Value = 98.77106105430903
$ws.Range("R8").Value = 888.9395494887812
$ws.Range("S8").Value = 0.05640701892459014
$ws.Range("T8").Value = 0.05640701892459014
$ws.Range("I9").Value = 0.1787346690539575
$ws.Range("J9").Value = 0.1787346690539575
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3186579999999999
$ws.Range("N9").Value = 0.9559739999999999
$ws.Range("O9").Value = 0.06962461334893082
$ws.Range("P9").Value = 0.06962461334893082
$ws.Range("Q9").Value = 21.790548789894
$ws.Range("R9").Value = 196.114939109046
$ws.Range("S9").Value = 0.0124443322249309
$ws.Range("T9").Value = 0.0124443322249309
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("M10").Value = 0.6327629999999999
$ws.Range("N10").Value = 1.898289
$ws.Range("O10").Value = 0.1382544270550543
$ws.Range("P10").Value = 0.1382544270550544
$ws.Range("Q10").Value = 33.71021723374199
$ws.Range("R10").Value = 303.3919551036779
$ws.Range("S10").Value = 0.01925151801710622
$ws.Range("T10").Value = 0.01925151801710622
$ws.Range("G11").Value = 53.27463399999999
$ws.Range("H11").Value = 159.823902
$ws.Range("I11").Value = 0.1392470275793777
$ws.Range("J11").Value = 0.1392470275793778
$ws.Range("O11").Value = 0.4765301499162115
$ws.Range("P11").Value = 0.4765301499162115
$ws.Range("Q11").Value = 116.1911066016446
$ws.Range("R11").Value = 1045.719959414802
$ws.Range("S11").Value = 0.0663554069277877
$ws.Range("T11").Value = 0.06635540692778773
$ws.Range("G12").Value = 53.27463399999999
$ws.Range("H12").Value = 159.823902
$ws.Range("I12").Value = 0.1392470275793777
$ws.Range("J12").Value = 0.1392470275793778
$ws.Range("M12").Value = 1.444396333333334
$ws.Range("N12").Value = 4.333189000000001
$ws.Range("O12").Value = 0.3155908096798033
$ws.Range("P12").Value = 0.3155908096798033
$ws.Range("Q12").Value = 76.94968600927534
$ws.Range("R12").Value = 692.547174083478
$ws.Range("S12").Value = 0.04394508217928172
$ws.Range("T12").Value = 0.04394508217928173
$ws.Range("G13").Value = 53.27463399999999
$ws.Range("H13").Value = 159.823902
$ws.Range("I13").Value = 0.1392470275793777
$ws.Range("J13").Value = 0.1392470275793778
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.3186579999999999
$ws.Range("N13").Value = 0.9559739999999999
$ws.Range("O13").Value = 0.06962461334893082
$ws.Range("P13").Value = 0.06962461334893082
$ws.Range("Q13").Value = 16.976388321172
$ws.Range("R13").Value = 152.7874948905479
$ws.Range("S13").Value = 0.009695020455202081
$ws.Range("T13").Value = 0.009695020455202083
$ws.Range("G14").Value = 91.01828266666666
$ws.Range("H14").Value = 273.054848
$ws.Range("I14").Value = 0.2378998101932138
$ws.Range("J14").Value = 0.2378998101932138
$ws.Range("M14").Value = 0.6327629999999999
$ws.Range("N14").Value = 1.898289
$ws.Range("O14").Value = 0.1382544270550543
$ws.Range("P14").Value = 0.1382544270550544
$ws.Range("Q14").Value = 57.59300159500798
$ws.Range("R14").Value = 518.3370143550719
$ws.Range("S14").Value = 0.03289070195476895
$ws.Range("T14").Value = 0.03289070195476896
$ws.Range("G15").Value = 91.01828266666666
$ws.Range("H15").Value = 273.054848
$ws.Range("I15").Value = 0.2378998101932138
$ws.Range("J15").Value = 0.2378998101932138
$ws.Range("O15").Value = 0.4765301499162115
$ws.Range("P15").Value = 0.4765301499162115
$ws.Range("Q15").Value = 198.5093878640498
$ws.Range("R15").Value = 1786.584490776448
$ws.Range("S15").Value = 0.1133664322164104
$ws.Range("T15").Value = 0.1133664322164104
$ws.Range("G16").Value = 91.01828266666666
$ws.Range("H16").Value = 273.054848
$ws.Range("I16").Value = 0.2378998101932138
$ws.Range("J16").Value = 0.2378998101932138
$ws.Range("M16").Value = 1.444396333333334
$ws.Range("N16").Value = 4.333189000000001
$ws.Range("O16").Value = 0.3155908096798033
$ws.Range("P16").Value = 0.3155908096798033
$ws.Range("Q16").Value = 131.4664737500302
$ws.Range("R16").Value = 1183.198263750272
$ws.Range("S16").Value = 0.07507899372154786
$ws.Range("T16").Value = 0.07507899372154787
$ws.Range("G17").Value = 91.01828266666666
$ws.Range("H17").Value = 273.054848
$ws.Range("I17").Value = 0.2378998101932138
$ws.Range("J17").Value = 0.2378998101932138
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.3186579999999999
$ws.Range("N17").Value = 0.9559739999999999
$ws.Range("O17").Value = 0.06962461334893082
$ws.Range("P17").Value = 0.06962461334893082
$ws.Range("Q17").Value = 29.00370391799466
$ws.Range("R17").Value = 261.0333352619519
$ws.Range("S17").Value = 0.01656368230048654
$ws.Range("T17").Value = 0.01656368230048654
